$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

# Passenger LDV lifetime: share new per year goes from 1/20 to 1/17
# (update the master formula cell B2 and each member of the shared
# formula group C2:H2 individually so the shared-formula grouping is
# preserved exactly as Excel originally laid it out).
$ws.Range("B2").Formula = "=1/17"
$ws.Range("C2").Formula = "=1/17"
$ws.Range("D2").Formula = "=1/17"
$ws.Range("E2").Formula = "=1/17"
$ws.Range("F2").Formula = "=1/17"
$ws.Range("G2").Formula = "=1/17"
$ws.Range("H2").Formula = "=1/17"

# Leave the cursor where the author left it when saving.
$ws.Range("F9").Select() | Out-Null
